$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.064.18"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "3.246.44"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'543.07"
$ws.Range("E5").Value = "  +2.65%  "
$ws.Range("D6").Value = "'147.72"
$ws.Range("E6").Value = "  +5.42%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "'7.38"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'0.115"
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "3.802.56"
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").Value = "'26.26"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "'0.0000175"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").Value = "61.012.56"
$ws.Range("E16").Value = "  +3.95%  "
$ws.Range("D17").Value = "3.218.41"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("D20").Value = "'8.38"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").Value = "'378.88"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'0.529"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'70.23"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("D26").Value = "'8.66"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +6.22%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.92"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.62"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'6.24"
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("E32").Value = "  +4.69%  "
$ws.Range("E33").Value = "  +6.80%  "
$ws.Range("D34").Value = "'6.63"
$ws.Range("E34").Value = "  +5.10%  "
$ws.Range("D35").Value = "'159.00"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("D37").Value = "'26.51"
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("D38").Value = "2.819.24"
$ws.Range("E38").Value = "  +4.50%  "
$ws.Range("D39").Value = "'0.0717"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("E40").Value = "  +8.74%  "
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").Value = "'4.29"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'39.98"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "3.284.62"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D47").Value = "'1.01"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "'21.13"
$ws.Range("E48").Value = "  +5.34%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "'0.809"
$ws.Range("E50").Value = "  +8.07%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.03%  "
